# Auto-generated Excel COM-interop script applying scheduled Kraken_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4177.5
$ws.Range("I116").Value = 3892.5
$ws.Range("J116").Value = 4747.5
$ws.Range("K116").Value = 3892.5
$ws.Range("L116").Value = 4747.5
$ws.Range("M116").Value = -450.5
$ws.Range("N116").Value = -11631.5

$ws.Range("H132").Value = 4771.273
$ws.Range("I132").Value = 3184.75
$ws.Range("J132").Value = 9002
$ws.Range("K132").Value = 9554.25
$ws.Range("L132").Value = 27006
$ws.Range("M132").Value = -7024.25
$ws.Range("N132").Value = -32066

$ws.Range("H141").Value = 5616.5264
$ws.Range("I141").Value = 5616.5264
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 16849.5792
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11669.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2400.6
$ws.Range("I2").Value = 2644
$ws.Range("J2").Value = 1832.6666
$ws.Range("K2").Value = 2644
$ws.Range("L2").Value = 1832.6666
$ws.Range("M2").Value = -2531
$ws.Range("N2").Value = -2058.6666

$ws.Range("H110").Value = 1764.2778
$ws.Range("I110").Value = 1133.8
$ws.Range("J110").Value = 4916.6665
$ws.Range("K110").Value = 1133.8
$ws.Range("L110").Value = 4916.6665
$ws.Range("M110").Value = 911.2
$ws.Range("N110").Value = -9006.666499999999

$ws.Range("H116").Value = 2400.6
$ws.Range("I116").Value = 2644
$ws.Range("J116").Value = 1832.6666
$ws.Range("K116").Value = 2644
$ws.Range("L116").Value = 1832.6666
$ws.Range("M116").Value = -350
$ws.Range("N116").Value = -6420.6666

$ws.Range("H132").Value = 1342.5385
$ws.Range("I132").Value = 1419.5834
$ws.Range("J132").Value = 418
$ws.Range("K132").Value = 4258.7502
$ws.Range("L132").Value = 1254
$ws.Range("M132").Value = -1728.7502
$ws.Range("N132").Value = -6314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2400.6
$ws.Range("I3").Value = 2644
$ws.Range("J3").Value = 1832.6666
$ws.Range("K3").Value = 2644
$ws.Range("L3").Value = 1832.6666
$ws.Range("M3").Value = -2530
$ws.Range("N3").Value = -2060.6666

$ws.Range("H99").Value = 3463.5454
$ws.Range("I99").Value = 3609.9
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 3609.9
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -2111.9

$ws.Range("H134").Value = 5289.1763
$ws.Range("I134").Value = 3608.75
$ws.Range("J134").Value = 9322.200000000001
$ws.Range("K134").Value = 10826.25
$ws.Range("L134").Value = 27966.6
$ws.Range("M134").Value = -8291.25
$ws.Range("N134").Value = -33036.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9108
$ws.Range("I31").Value = 9274.375
$ws.Range("J31").Value = 8664.333000000001
$ws.Range("K31").Value = 9274.375
$ws.Range("L31").Value = 8664.333000000001
$ws.Range("M31").Value = -8979.375
$ws.Range("N31").Value = -9254.333000000001

$ws.Range("H34").Value = 9108
$ws.Range("I34").Value = 9274.375
$ws.Range("J34").Value = 8664.333000000001
$ws.Range("K34").Value = 9274.375
$ws.Range("L34").Value = 8664.333000000001
$ws.Range("M34").Value = -9072.375
$ws.Range("N34").Value = -9068.333000000001

$ws.Range("H58").Value = 2971.25
$ws.Range("I58").Value = 3211.1538
$ws.Range("J58").Value = 1931.6666
$ws.Range("K58").Value = 3211.1538
$ws.Range("L58").Value = 1931.6666
$ws.Range("M58").Value = -3008.1538
$ws.Range("N58").Value = -2337.6666

$ws.Range("H105").Value = 1253.3334
$ws.Range("I105").Value = 1253.3334
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1253.3334
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 493.6666
$ws.Range("N105").ClearContents()

$ws.Range("H134").Value = 2010
$ws.Range("I134").Value = 1893.8096
$ws.Range("J134").Value = 2620
$ws.Range("K134").Value = 5681.4288
$ws.Range("L134").Value = 7860
$ws.Range("M134").Value = -3146.4288
$ws.Range("N134").Value = -12930

$ws.Range("H136").Value = 2971.25
$ws.Range("I136").Value = 3211.1538
$ws.Range("J136").Value = 1931.6666
$ws.Range("K136").Value = 9633.4614
$ws.Range("L136").Value = 5794.9998
$ws.Range("M136").Value = -7083.4614
$ws.Range("N136").Value = -10894.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 150
$ws.Range("I25").Value = 150
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 450
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -281

$ws.Range("H30").Value = 150
$ws.Range("I30").Value = 150
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 450
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -348

$ws.Range("H124").Value = 8000
$ws.Range("I124").Value = 8000
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 24000
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = -19090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 30000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 30000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 30000
$ws.Range("N33").Value = -30504

$ws.Range("H44").Value = 25000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 25000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -26192
$ws.Range("M44").ClearContents()

$ws.Range("H62").Value = 90000
$ws.Range("I62").Value = 90000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 90000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -89314

$ws.Range("H65").Value = 90000
$ws.Range("I65").Value = 90000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 270000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -266568

$ws.Range("H102").Value = 3059.5
$ws.Range("I102").Value = 2877.4443
$ws.Range("J102").Value = 4698
$ws.Range("K102").Value = 2877.4443
$ws.Range("L102").Value = 4698
$ws.Range("M102").Value = -1255.4443

$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1370

$ws.Range("H132").Value = 2330.3333
$ws.Range("I132").Value = 1965.9412
$ws.Range("J132").Value = 3879
$ws.Range("K132").Value = 5897.8236
$ws.Range("L132").Value = 11637
$ws.Range("M132").Value = -3367.8236
$ws.Range("N132").Value = -16697

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5197.154
$ws.Range("I7").Value = 5233.1816
$ws.Range("J7").Value = 4999
$ws.Range("K7").Value = 5233.1816
$ws.Range("L7").Value = 4999
$ws.Range("M7").Value = -5121.1816
$ws.Range("N7").Value = -5223

$ws.Range("H22").Value = 1647.5
$ws.Range("I22").Value = 1530
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1530
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1235
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 1647.5
$ws.Range("I27").Value = 1530
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1530
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1423
$ws.Range("N27").Value = -2214

$ws.Range("H40").Value = 4248.0835
$ws.Range("I40").Value = 4248.0835
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 4248.0835
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -4112.0835

$ws.Range("H93").Value = 3134.8
$ws.Range("I93").Value = 3762.1667
$ws.Range("J93").Value = 2193.75
$ws.Range("K93").Value = 3762.1667
$ws.Range("L93").Value = 2193.75
$ws.Range("M93").Value = -2514.1667
$ws.Range("N93").Value = -4689.75

$ws.Range("H100").Value = 6626.3076
$ws.Range("I100").Value = 3186
$ws.Range("J100").Value = 12130.8
$ws.Range("K100").Value = 3186
$ws.Range("L100").Value = 12130.8
$ws.Range("M100").Value = -2645
$ws.Range("N100").Value = -13212.8

$ws.Range("H122").Value = 3260
$ws.Range("I122").Value = 3260
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9780
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7330
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 5197.154
$ws.Range("I126").Value = 5233.1816
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 15699.5448
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -13229.5448
$ws.Range("N126").Value = -19937

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 6250
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 18750
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -16220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1631.8182
$ws.Range("I136").Value = 1781.25
$ws.Range("J136").Value = 1233.3334
$ws.Range("K136").Value = 5343.75
$ws.Range("L136").Value = 3700.0002
$ws.Range("M136").Value = -2793.75
$ws.Range("N136").Value = -8800.0002

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
